# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-22 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 4
    4  = 2
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 0
    10 = 3
    11 = 1
    12 = 1
    13 = 2
    14 = 3
    15 = 1
    16 = 7
    17 = 2
    18 = 3
    19 = 4
    20 = 1
    21 = 2
    22 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
